$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $para = $cell.Range.Paragraphs.Item(1).Range
    $contentRange = $d.Range($para.Start, $para.End - 1)
    $contentRange.Text = $newText
}

Set-CellText $t 1 1 "15÷7="
Set-CellText $t 1 2 "35÷7="
Set-CellText $t 1 3 "58÷4="
Set-CellText $t 1 4 "95÷3="
Set-CellText $t 1 5 "78÷4="
Set-CellText $t 5 1 "70÷4="
Set-CellText $t 5 2 "42÷3="
Set-CellText $t 5 3 "33÷4="
Set-CellText $t 5 4 "21÷5="
Set-CellText $t 5 5 "98÷5="
Set-CellText $t 9 1 "24÷3="
Set-CellText $t 9 2 "80÷3="
Set-CellText $t 9 3 "52÷3="
Set-CellText $t 9 4 "47÷9="
Set-CellText $t 9 5 "96÷9="
Set-CellText $t 13 1 "15÷7="
Set-CellText $t 13 2 "12÷8="
Set-CellText $t 13 3 "34÷8="
Set-CellText $t 13 4 "80÷3="
Set-CellText $t 13 5 "90÷7="
Set-CellText $t 17 1 "63÷5="
Set-CellText $t 17 2 "81÷4="
Set-CellText $t 17 3 "26÷6="
Set-CellText $t 17 4 "40÷9="
Set-CellText $t 17 5 "69÷2="

Write-Host "Replacements applied."
